$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 2020 column (K) of data, mirroring the style of the corresponding
# 2019 column (J) cell in each row, then clear the inherited vertical-center
# alignment so the new style entries stay as plain font+border combos.

$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("K4").Value = 2020
$ws.Range("K4").VerticalAlignment = -4107

$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("K5").Value = 22
$ws.Range("K5").VerticalAlignment = -4107

$ws.Range("J6").Copy()
$ws.Range("K6").PasteSpecial(-4122)
$ws.Range("K6").Value = 29.1
$ws.Range("K6").VerticalAlignment = -4107

$ws.Range("J7").Copy()
$ws.Range("K7").PasteSpecial(-4122)
$ws.Range("K7").Value = 20.2
$ws.Range("K7").VerticalAlignment = -4107

$ws.Range("J8").Copy()
$ws.Range("K8").PasteSpecial(-4122)
$ws.Range("K8").Value = 26.8
$ws.Range("K8").VerticalAlignment = -4107

$ws.Range("J9").Copy()
$ws.Range("K9").PasteSpecial(-4122)
$ws.Range("K9").Value = 39.8
$ws.Range("K9").VerticalAlignment = -4107

$ws.Range("J10").Copy()
$ws.Range("K10").PasteSpecial(-4122)
$ws.Range("K10").Value = 22.7
$ws.Range("K10").VerticalAlignment = -4107

$ws.Range("J11").Copy()
$ws.Range("K11").PasteSpecial(-4122)
$ws.Range("K11").Value = 22
$ws.Range("K11").VerticalAlignment = -4107

$ws.Range("J12").Copy()
$ws.Range("K12").PasteSpecial(-4122)
$ws.Range("K12").Value = 33.8
$ws.Range("K12").VerticalAlignment = -4107

$ws.Range("J13").Copy()
$ws.Range("K13").PasteSpecial(-4122)
$ws.Range("K13").Value = 18.8
$ws.Range("K13").VerticalAlignment = -4107

$ws.Range("J14").Copy()
$ws.Range("K14").PasteSpecial(-4122)
$ws.Range("K14").Value = 28
$ws.Range("K14").VerticalAlignment = -4107

$ws.Range("J15").Copy()
$ws.Range("K15").PasteSpecial(-4122)
$ws.Range("K15").Value = 38.7
$ws.Range("K15").VerticalAlignment = -4107

$ws.Range("J16").Copy()
$ws.Range("K16").PasteSpecial(-4122)
$ws.Range("K16").Value = 22
$ws.Range("K16").VerticalAlignment = -4107

$ws.Range("J17").Copy()
$ws.Range("K17").PasteSpecial(-4122)
$ws.Range("K17").Value = 35.1
$ws.Range("K17").VerticalAlignment = -4107

$ws.Range("J18").Copy()
$ws.Range("K18").PasteSpecial(-4122)
$ws.Range("K18").Value = 45.1
$ws.Range("K18").VerticalAlignment = -4107

$ws.Range("J19").Copy()
$ws.Range("K19").PasteSpecial(-4122)
$ws.Range("K19").Value = 33.3
$ws.Range("K19").VerticalAlignment = -4107

$ws.Range("J20").Copy()
$ws.Range("K20").PasteSpecial(-4122)
$ws.Range("K20").Value = 19.4
$ws.Range("K20").VerticalAlignment = -4107

$ws.Range("J21").Copy()
$ws.Range("K21").PasteSpecial(-4122)
$ws.Range("K21").Value = 13
$ws.Range("K21").VerticalAlignment = -4107

$ws.Range("J22").Copy()
$ws.Range("K22").PasteSpecial(-4122)
$ws.Range("K22").Value = 19.9
$ws.Range("K22").VerticalAlignment = -4107

$ws.Range("J23").Copy()
$ws.Range("K23").PasteSpecial(-4122)
$ws.Range("K23").Value = 26.2
$ws.Range("K23").VerticalAlignment = -4107

$ws.Range("J24").Copy()
$ws.Range("K24").PasteSpecial(-4122)
$ws.Range("K24").Value = 52.8
$ws.Range("K24").VerticalAlignment = -4107

$ws.Range("J25").Copy()
$ws.Range("K25").PasteSpecial(-4122)
$ws.Range("K25").Value = 22.5
$ws.Range("K25").VerticalAlignment = -4107

$ws.Range("J26").Copy()
$ws.Range("K26").PasteSpecial(-4122)
$ws.Range("K26").Value = 20.1
$ws.Range("K26").VerticalAlignment = -4107

$ws.Range("J27").Copy()
$ws.Range("K27").PasteSpecial(-4122)
$ws.Range("K27").Value = 33.6
$ws.Range("K27").VerticalAlignment = -4107

$ws.Range("J28").Copy()
$ws.Range("K28").PasteSpecial(-4122)
$ws.Range("K28").Value = 16.9
$ws.Range("K28").VerticalAlignment = -4107

$ws.Range("J29").Copy()
$ws.Range("K29").PasteSpecial(-4122)
$ws.Range("K29").Value = 23.5
$ws.Range("K29").VerticalAlignment = -4107

$ws.Range("J30").Copy()
$ws.Range("K30").PasteSpecial(-4122)
$ws.Range("K30").Value = 30.7
$ws.Range("K30").VerticalAlignment = -4107

$excel.CutCopyMode = $false

# Match the authors final selection recorded in the workbook view.
$ws.Range("K18").Select()
